# Update the "number of attendees/favorites" column F figures on the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets to the freshly scraped
# figures (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 794
$ws1.Range("F10").Value = 274
$ws1.Range("F12").Value = 10300
$ws1.Range("F16").Value = 647
$ws1.Range("F17").Value = 11866
$ws1.Range("F18").Value = 12253
$ws1.Range("F21").Value = 17

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 794
$ws4.Range("F11").Value = 274
$ws4.Range("F13").Value = 10300
$ws4.Range("F17").Value = 647
$ws4.Range("F18").Value = 11866
$ws4.Range("F19").Value = 12253
$ws4.Range("F22").Value = 17
